# Auto-generated edit script for COLORADO_2015-style workbook cleanup
# - Renames header columns to short codes
# - Normalizes connector words (de/del/la/las/el/los/y) to capitalized form
#   in municipality/state names (mirrors an updated cleaning script)
# - Removes trailing metadata/footer rows 1378-1382

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections (header renames + capitalization fixes) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B11").Value = "San José De Gracia"
$ws.Range("B16").Value = "Playas De Rosarito"
$ws.Range("B33").Value = "Amatenango De La Frontera"
$ws.Range("B38").Value = "Chiapa De Corzo"
$ws.Range("B55").Value = "Mazapa De Madero"
$ws.Range("B63").Value = "Salto De Agua"
$ws.Range("B64").Value = "San Cristóbal De Las Casas"
$ws.Range("B102").Value = "Coyame Del Sotol"
$ws.Range("B113").Value = "Guadalupe Y Calvo"
$ws.Range("B116").Value = "Hidalgo Del Parral"
$ws.Range("B142").Value = "San Francisco De Borja"
$ws.Range("B143").Value = "San Francisco De Conchos"
$ws.Range("B144").Value = "San Francisco Del Oro"
$ws.Range("B152").Value = "Valle De Zaragoza"
$ws.Range("B167").Value = "San Juan De Sabinas"
$ws.Range("B181").Value = "Villa De Álvarez"
$ws.Range("A183").Value = "Ciudad De México"
$ws.Range("B187").Value = "Cuajimalpa De Morelos"
$ws.Range("B202").Value = "Coneto De Comonfort"
$ws.Range("B216").Value = "Nombre De Dios"
$ws.Range("B220").Value = "Pánuco De Coronado"
$ws.Range("B227").Value = "San Juan De Guadalupe"
$ws.Range("B228").Value = "San Juan Del Río"
$ws.Range("B229").Value = "San Luis Del Cordero"
$ws.Range("B230").Value = "San Pedro Del Gallo"
$ws.Range("A240").Value = "Estado De México"
$ws.Range("B240").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B242").Value = "Almoloya De Alquisiras"
$ws.Range("B243").Value = "Almoloya De Juárez"
$ws.Range("B248").Value = "Atizapán De Zaragoza"
$ws.Range("B254").Value = "Coacalco De Berriozábal"
$ws.Range("B260").Value = "Ecatepec De Morelos"
$ws.Range("B265").Value = "Ixtapan De La Sal"
$ws.Range("B266").Value = "Ixtapan Del Oro"
$ws.Range("B277").Value = "Naucalpan De Juárez"
$ws.Range("B285").Value = "San Felipe Del Progreso"
$ws.Range("B286").Value = "San Martín De Las Pirámides"
$ws.Range("B288").Value = "San Simón De Guerrero"
$ws.Range("B290").Value = "Soyaniquilpan De Juárez"
$ws.Range("B298").Value = "Tenango Del Valle"
$ws.Range("B306").Value = "Tlalnepantla De Baz"
$ws.Range("B311").Value = "Valle De Bravo"
$ws.Range("B312").Value = "Villa De Allende"
$ws.Range("B313").Value = "Villa Del Carbón"
$ws.Range("B323").Value = "San Miguel De Allende"
$ws.Range("B324").Value = "Apaseo El Alto"
$ws.Range("B325").Value = "Apaseo El Grande"
$ws.Range("B333").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B337").Value = "Jaral Del Progreso"
$ws.Range("B345").Value = "Purísima Del Rincón"
$ws.Range("B350").Value = "San Francisco Del Rincón"
$ws.Range("B352").Value = "San Luis De La Paz"
$ws.Range("B354").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B355").Value = "Silao De La Victoria"
$ws.Range("B359").Value = "Valle De Santiago"
$ws.Range("B365").Value = "Acapulco De Juárez"
$ws.Range("B367").Value = "Ajuchitlán Del Progreso"
$ws.Range("B368").Value = "Alcozauca De Guerrero"
$ws.Range("B371").Value = "Atenango Del Río"
$ws.Range("B372").Value = "Atoyac De Álvarez"
$ws.Range("B373").Value = "Ayutla De Los Libres"
$ws.Range("B375").Value = "Buenavista De Cuéllar"
$ws.Range("B376").Value = "Chilapa De Álvarez"
$ws.Range("B377").Value = "Chilpancingo De Los Bravo"
$ws.Range("B378").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B383").Value = "Coyuca De Benítez"
$ws.Range("B384").Value = "Coyuca De Catalán"
$ws.Range("B388").Value = "Cuetzala Del Progreso"
$ws.Range("B389").Value = "Cutzamala De Pinzón"
$ws.Range("B395").Value = "Huitzuco De Los Figueroa"
$ws.Range("B396").Value = "Iguala De La Independencia"
$ws.Range("B398").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B401").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B404").Value = "Mártir De Cuilapan"
$ws.Range("B417").Value = "Taxco De Alarcón"
$ws.Range("B419").Value = "Técpan De Galeana"
$ws.Range("B421").Value = "Tepecoacuilco De Trujano"
$ws.Range("B423").Value = "Tixtla De Guerrero"
$ws.Range("B426").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B427").Value = "Tlapa De Comonfort"
$ws.Range("B442").Value = "Atotonilco El Grande"
$ws.Range("B446").Value = "Cuautepec De Hinojosa"
$ws.Range("B450").Value = "Huasca De Ocampo"
$ws.Range("B452").Value = "Huejutla De Reyes"
$ws.Range("B455").Value = "Jacala De Ledezma"
$ws.Range("B460").Value = "Mineral Del Chico"
$ws.Range("B461").Value = "Mineral Del Monte"
$ws.Range("B462").Value = "Mixquiahuala De Juárez"
$ws.Range("B464").Value = "Pachuca De Soto"
$ws.Range("B465").Value = "Progreso De Obregón"
$ws.Range("B470").Value = "Tenango De Doria"
$ws.Range("B472").Value = "Tepehuacán De Guerrero"
$ws.Range("B473").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B475").Value = "Tezontepec De Aldama"
$ws.Range("B480").Value = "Tula De Allende"
$ws.Range("B481").Value = "Tulancingo De Bravo"
$ws.Range("B482").Value = "Zapotlán De Juárez"
$ws.Range("B487").Value = "Acatlán De Juárez"
$ws.Range("B488").Value = "Ahualulco De Mercado"
$ws.Range("B491").Value = "Atemajac De Brizuela"
$ws.Range("B494").Value = "Atotonilco El Alto"
$ws.Range("B496").Value = "Autlán De Navarro"
$ws.Range("B502").Value = "Cañadas De Obregón"
$ws.Range("B507").Value = "Cuautitlán De García Barragán"
$ws.Range("B515").Value = "Encarnación De Díaz"
$ws.Range("B521").Value = "Huejuquilla El Alto"
$ws.Range("B522").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B523").Value = "Ixtlahuacán Del Río"
$ws.Range("B527").Value = "Jilotlán De Los Dolores"
$ws.Range("B533").Value = "Lagos De Moreno"
$ws.Range("B539").Value = "Ojuelos De Jalisco"
$ws.Range("B544").Value = "San Diego De Alejandría"
$ws.Range("B545").Value = "San Juan De Los Lagos"
$ws.Range("B546").Value = "San Juanito De Escobedo"
$ws.Range("B548").Value = "San Martín De Bolaños"
$ws.Range("B550").Value = "San Miguel El Alto"
$ws.Range("B551").Value = "San Sebastián Del Oeste"
$ws.Range("B552").Value = "Santa María De Los Ángeles"
$ws.Range("B555").Value = "Talpa De Allende"
$ws.Range("B556").Value = "Tamazula De Gordiano"
$ws.Range("B560").Value = "Teocuitatlán De Corona"
$ws.Range("B561").Value = "Tepatitlán De Morelos"
$ws.Range("B563").Value = "Tizapán El Alto"
$ws.Range("B564").Value = "Tlajomulco De Zúñiga"
$ws.Range("B576").Value = "Unión De San Antonio"
$ws.Range("B577").Value = "Unión De Tula"
$ws.Range("B578").Value = "Valle De Guadalupe"
$ws.Range("B579").Value = "Valle De Juárez"
$ws.Range("B584").Value = "Yahualica De González Gallo"
$ws.Range("B585").Value = "Zacoalco De Torres"
$ws.Range("B588").Value = "Zapotlán Del Rey"
$ws.Range("B589").Value = "Zapotlán El Grande"
$ws.Range("B610").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B671").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B695").Value = "Coatlán Del Río"
$ws.Range("B705").Value = "Puente De Ixtla"
$ws.Range("B711").Value = "Tetela Del Volcán"
$ws.Range("B712").Value = "Tlaltizapán De Zapata"
$ws.Range("B717").Value = "Zacualpan De Amilpas"
$ws.Range("B721").Value = "Amatlán De Cañas"
$ws.Range("B722").Value = "Bahía De Banderas"
$ws.Range("B726").Value = "Ixtlán Del Río"
$ws.Range("B733").Value = "Santa María Del Oro"
$ws.Range("B742").Value = "Ciénega De Flores"
$ws.Range("B751").Value = "Mier Y Noriega"
$ws.Range("B755").Value = "San Nicolás De Los Garza"
$ws.Range("B758").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B761").Value = "Ayoquezco De Aldama"
$ws.Range("B764").Value = "Chalcatongo De Hidalgo"
$ws.Range("B765").Value = "Coicoyán De Las Flores"
$ws.Range("B768").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B769").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B771").Value = "Ixtlán De Juárez"
$ws.Range("B772").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B778").Value = "Mariscala De Juárez"
$ws.Range("B780").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B782").Value = "Nejapa De Madero"
$ws.Range("B783").Value = "Oaxaca De Juárez"
$ws.Range("B784").Value = "Ocotlán De Morelos"
$ws.Range("B786").Value = "Putla Villa De Guerrero"
$ws.Range("B790").Value = "San Antonio De La Cal"
$ws.Range("B808").Value = "San Juan Del Estado"
$ws.Range("B828").Value = "San Miguel Del Puerto"
$ws.Range("B832").Value = "San Pablo Villa De Mitla"
$ws.Range("B843").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B855").Value = "Santa María Del Tule"
$ws.Range("B865").Value = "Santiago Del Río"
$ws.Range("B885").Value = "Sitio De Xitlapehua"
$ws.Range("B886").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B888").Value = "Teotitlán De Flores Magón"
$ws.Range("B889").Value = "Teotitlán Del Valle"
$ws.Range("B890").Value = "Tlacolula De Matamoros"
$ws.Range("B892").Value = "Totontepec Villa De Morelos"
$ws.Range("B894").Value = "Villa De Chilapa De Díaz"
$ws.Range("B895").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B896").Value = "Villa De Zaachila"
$ws.Range("B899").Value = "Villa Sola De Vega"
$ws.Range("B900").Value = "Zimatlán De Álvarez"
$ws.Range("B914").Value = "Chalchicomula De Sesma"
$ws.Range("B925").Value = "Cuayuca De Andrade"
$ws.Range("B932").Value = "Huehuetlán El Chico"
$ws.Range("B933").Value = "Huehuetlán El Grande"
$ws.Range("B937").Value = "Ixcamilpa De Guerrero"
$ws.Range("B940").Value = "Izúcar De Matamoros"
$ws.Range("B951").Value = "Palmar De Bravo"
$ws.Range("B961").Value = "San Nicolás De Los Ranchos"
$ws.Range("B964").Value = "San Salvador El Seco"
$ws.Range("B965").Value = "San Salvador El Verde"
$ws.Range("B969").Value = "Tecali De Herrera"
$ws.Range("B975").Value = "Tepanco De López"
$ws.Range("B976").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B980").Value = "Tepexi De Rodríguez"
$ws.Range("B982").Value = "Tetela De Ocampo"
$ws.Range("B987").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B996").Value = "Totoltepec De Guerrero"
$ws.Range("B1009").Value = "Amealco De Bonfil"
$ws.Range("B1011").Value = "Cadereyta De Montes"
$ws.Range("B1015").Value = "Jalpan De Serra"
$ws.Range("B1016").Value = "Landa De Matamoros"
$ws.Range("B1018").Value = "Pinal De Amoles"
$ws.Range("B1021").Value = "San Juan Del Río"
$ws.Range("B1031").Value = "Ciudad Del Maíz"
$ws.Range("B1039").Value = "Mexquitic De Carmona"
$ws.Range("B1044").Value = "San Ciro De Acosta"
$ws.Range("B1049").Value = "Santa María Del Río"
$ws.Range("B1051").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1055").Value = "Tanquián De Escobedo"
$ws.Range("B1058").Value = "Villa De Arriaga"
$ws.Range("B1059").Value = "Villa De Guadalupe"
$ws.Range("B1060").Value = "Villa De Ramos"
$ws.Range("B1061").Value = "Villa De Reyes"
$ws.Range("B1106").Value = "Nacozari De García"
$ws.Range("B1113").Value = "San Miguel De Horcasitas"
$ws.Range("B1124").Value = "Jalpa De Méndez"
$ws.Range("B1148").Value = "Soto La Marina"
$ws.Range("B1161").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1164").Value = "Muñoz De Domingo Arenas"
$ws.Range("B1165").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1167").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1171").Value = "Tepetitla De Lardizábal"
$ws.Range("B1189").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1191").Value = "Amatlán De Los Reyes"
$ws.Range("B1197").Value = "Castillo De Teayo"
$ws.Range("B1199").Value = "Cazones De Herrera"
$ws.Range("B1210").Value = "Cosamaloapan De Carpio"
$ws.Range("B1224").Value = "Hueyapan De Ocampo"
$ws.Range("B1225").Value = "Ignacio De La Llave"
$ws.Range("B1228").Value = "Ixhuatlán De Madero"
$ws.Range("B1229").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1237").Value = "Juchique De Ferrer"
$ws.Range("B1240").Value = "Las Vigas De Ramírez"
$ws.Range("B1241").Value = "Lerdo De Tejada"
$ws.Range("B1244").Value = "Martínez De La Torre"
$ws.Range("B1246").Value = "Medellín De Bravo"
$ws.Range("B1260").Value = "Paso De Ovejas"
$ws.Range("B1261").Value = "Paso Del Macho"
$ws.Range("B1265").Value = "Poza Rica De Hidalgo"
$ws.Range("B1274").Value = "Soledad De Doblado"
$ws.Range("B1293").Value = "Vega De Alatorre"
$ws.Range("B1300").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1328").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1330").Value = "Concepción Del Oro"
$ws.Range("B1341").Value = "Jiménez Del Teul"
$ws.Range("B1348").Value = "Mezquital Del Oro"
$ws.Range("B1352").Value = "Nochistlán De Mejía"
$ws.Range("B1353").Value = "Noria De Ángeles"
$ws.Range("B1364").Value = "Teúl De González Ortega"
$ws.Range("B1365").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1366").Value = "Trinidad García De La Cadena"
$ws.Range("B1369").Value = "Villa De Cos"

# --- Floating point recompute artifacts (count/total*100/100) ---
$ws.Range("D86").Value = 19/20238*100/100
$ws.Range("D352").Value = 19/20238*100/100
$ws.Range("D363").Value = 19/20238*100/100
$ws.Range("D374").Value = 19/20238*100/100
$ws.Range("D417").Value = 19/20238*100/100
$ws.Range("D650").Value = 19/20238*100/100
$ws.Range("D671").Value = 19/20238*100/100
$ws.Range("D731").Value = 19/20238*100/100
$ws.Range("D1149").Value = 19/20238*100/100
$ws.Range("D1335").Value = 19/20238*100/100

# --- Remove trailing metadata/footer rows ---
$ws.Range("A1378:D1382").EntireRow.Delete()

